$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.7598989048020219
$ws.Range("C2").Value = 0.7755803955288049
$ws.Range("D2").Value = 0.7676595744680852
$ws.Range("E2").Value = 1163

# Row 3
$ws.Range("B3").Value = 0.8288431061806656
$ws.Range("C3").Value = 0.812111801242236
$ws.Range("D3").Value = 0.8203921568627451
$ws.Range("E3").Value = 644

# Row 4
$ws.Range("B4").Value = 0.742483660130719
$ws.Range("C4").Value = 0.7319587628865979
$ws.Range("D4").Value = 0.7371836469824788
$ws.Range("E4").Value = 776

# Row 5
$ws.Range("B5").Value = 0.5154061624649859
$ws.Range("C5").Value = 0.5154061624649859
$ws.Range("D5").Value = 0.5154061624649859
$ws.Range("E5").Value = 357

# Row 6
$ws.Range("B6").Value = 0.7404761904761905
$ws.Range("C6").Value = 0.7404761904761905
$ws.Range("D6").Value = 0.7404761904761905
$ws.Range("E6").Value = 0.7404761904761905

# Row 7
$ws.Range("B7").Value = 0.7116579583945981
$ws.Range("C7").Value = 0.7087642805306562
$ws.Range("D7").Value = 0.7101603851945737
$ws.Range("E7").Value = 2940

# Row 8
$ws.Range("B8").Value = 0.7407158867097069
$ws.Range("C8").Value = 0.7404761904761905
$ws.Range("D8").Value = 0.7405357633280253
$ws.Range("E8").Value = 2940
